$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: date value, copy number format (date style) from A41
$ws.Range("A42").Value = 45890
$ws.Range("A42").NumberFormat = $ws.Range("A41").NumberFormat()

# Numeric OHLCV + return/label columns
$ws.Range("B42").Value = 6370.169921875
$ws.Range("C42").Value = 6393.64990234375
$ws.Range("D42").Value = 6352.7099609375
$ws.Range("E42").Value = 6380.830078125
$ws.Range("F42").Value = 4013560000
$ws.Range("G42").Value = -0.0040041815293089
$ws.Range("H42").Value = -1

# Headline text columns (identical text in I and J)
$headline = 'bcs veterinary college say it risk insolvency in if fee dont increase big mac price cut signal mcdonalds return to value a qsr battle intensify target ceo to step down amid company struggle target ceo to step down amid company struggle google pixel phone raise the ante on artificial intelligence qualified small business stock in tax and estate planning bitcoin could go to thanks to global retirement account could drive say bill miller iv instagrams new feature could change how you make friend online mastercard and alloy team to launch fraud prevention solution galbot clinch gold medal at world humanoid robot game redefining autonomous robotics in realworld application this is finally the end the ripple case concludes a predicted superior energy service announces acquisition of quail tool from nabors nearly doubling tubular inventory and expanding global reach canada seek meeting to address china wto complaint over steel tariff what is going on with micron stock wednesday nancy pelosi cathie wood make million on databricks new billion valuation inventhelp inventor develops new tool for removing tumbleweed mho cf financial corporation announces quarterly dividend ethereum whale with b portfolio back new ethbased meme coin next crypto to explode xrp news today xrp price drop below a whale selling and market pressure mount cf financial corporation announces quarterly dividend tech selloff push down nasdaq sp investor cautious ahead of fed meeting reuters solana price prediction stunning rally expected but it this new sol rival in the spotlight crypto news today yuanbacked stablecoins krakens xstocks on tron and pumpfuns memecoin surge a marvell stock price drop this bearish trade make sense brave dog who battled a bull earns hambone award for most unusual pet story will vibe coding spark a new economy elon musk courtroom drama take a turn she thought that her computer science degree would get her a six figure job instead it got her an interview with chipotle deadline approaching cto realty growth inc cto investor who lost money urged to contact law office of howard g smith fomc minute show fed member expect higher inflation the rise and fall of target fortune trump chilling effect is coming for museum historian warn bookingcom parent agrees to million settlement in texas hotel fee lawsuit whats going on with doordash stock wednesday inventhelp inventor develops new ice and snow removal tool for vehicle mbq why these camera are the best for photography hidden cocktail lounge to open above sd oyster co in uptown nabors announces sale of quail tool to superior energy service for million nabors announces sale of quail tool to superior energy service for million visa make online payment tech the star of it asia move ti chief financial officer rafael lizardi to speak at citi investor conference lovelands only hour mcdonalds now open spectrum pharmaceutical sppi security lawsuit new lead plaintiff deadline announced parent of girl killed in texas camp flooding say common sense safety measure were absent a lawmaker consider bill best wireless earbuds for workout under r august presale pick that could create millionaireswhy ozak ai is the most hyped ai token of perplexity b moonshot buying google chrome r lakh digital arrest scam bengaluru techie duped realpage debut ai platform with openai doordash add tiered membership package to commerce platform for restaurant fed dissenter appeared alone in favoring rate cut at july meeting minute show reuters summer extends wellbeyond labor day on miami beach with experience offering celebration of design wellness and culinary excellence talen energy stock a deep dive into analyst perspective rating where genpact stand with analyst how to trade spy and top tech stock using technical analysis where weatherford international stand with analyst divided fed worried about tariff inflation and the labor market minute show cnbc breaking down gaming and leisure prop analyst share their view what analyst rating have to say about gulfport energy bernie sander back trump plan to buy stake in intel ar technica how is the market feeling about ametek the best toaster oven of tried and tested looking into essex property trust recent short interest a look into bausch health company inc price over earnings is the market bullish or bearish on commvault system sofi technology stock a deep dive into analyst perspective rating americorps is under siege what happens in the community it serf inventhelp inventor develop wet extraction dna collection devicekit mbq hawkish fomc minute show majority fear higher inflation more than lower employment dogecoin news today dogecoin price action remains bullish supported by longterm indicator cold wallet raise m while hedera eye etf and chainlink target which top trending crypto win surf bio announces peerreviewed publication demonstrating feasibility of subcutaneous delivery of ultraconcentrated biologics nickel glut weighs on stainless steel despite u tariff support daybreaker and flavcity are throwing the biggest sober dance party in u history black hill northwestern energy to merge in billion deal from presale star to market high why bdag eth link doge are the top crypto to buy this august fuel station packed in russian far east a drone strike cause shortage housing advocate want more buyin from the state viking premium valuation backed by strong growth analyst note highquality eyeglass at great price is zennis mission life molecular imaging expands neuraceq florbetaben f injection production in southern california essentia health st joseph medical center brainerd birth give your wardrobe a refresh with champion off sale what will the prime minister call for a selfreliant india mean for new delhi trade and industrial policy vance derides protester while visiting national guard at union station bed bath beyond bos take a shot at california refusing to reopen store in the state san francisco chronicle deadline approaching kindercare learning company inc klc investor who lost money urged to contact law office of howard g smith get tracking with up to off the samsung galaxy ring united delta sued for allegedly selling windowless window seat wicys make stride in building the cybersecurity workforce with th annual security training scholarship calpaks labor day deal are here take up to off travel musthaves suspect in american idol exec murder case may be too mentally ill for prosecution ming shing group holding limited announces entering into a bitcoin purchase agreement for the purchase of bitcoins mcdonalds cut combo price to regain affordability edge report james hardie stock is trading lower wednesday whats going on target digital sale gain a traffic wane and ceo departs how to create highquality informative content top crypto presales of breakthrough contender with x potential billionaire backing woke math doesnt add up amid dei rollback court reject doj request to unseal jeffrey epstein grand jury file white house start tiktok account ahead of apps looming ban white house start tiktok account ahead of apps looming ban obama applauds newsoms california redistricting plan a responsible target ceo to step down amid company struggle obama applauds newsoms california redistricting plan a responsible a texas gop push new map obama applauds newsoms california redistricting plan a responsible a texas gop push new map google pixel phone raise the ante on artificial intelligence judge decline to unseal grand jury material in jeffrey epstein case a responsible approach obama back california response to texas new congressional line judge rule that some texas school dont have to display ten commandment in classroom both party expect a gop map in texas to clear a big hurdle in a national fight over redistricting russia say must be part of ukraine security guarantee talk here come the food desert kroger to close store in washington state thanks to high crime do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast french diplomacy macron just called putin an ogre political vampire are draining the country term limit are the answer do kid turn to phone because parent restrict their independence the ethical life podcast do kid turn to phone because parent restrict their independence the ethical life podcast epstein file latest grand jury transcript request denied what trump is really up to with the military occupation of dc usled coalition capture a senior islamic state member in syria the rise and fall of target fortune trump chilling effect is coming for museum historian warn federal judge denies trump admin request to unseal epstein grand jury material federal judge denies trump admin request to unseal epstein grand jury material federal judge denies trump admin request to unseal epstein grand jury material israel launch diplomatic attack on it western ally ahead of palestinian statehood recognition federal judge denies trump admin request to unseal epstein grand jury material texas cant put the ten commandment in certain school district classroom judge say texas cant put the ten commandment in certain school district classroom judge say full of s new york republican accuses state dems of hypocrisy in redistricting push judge decline to unseal grand jury material in jeffrey epstein case trump overture to putin gain nothing cost ukrainian dearly jonah goldberg trump overture to putin gain nothing cost ukrainian dearly jonah goldberg trump overture to putin gain nothing cost ukrainian dearly jonah goldberg trump overture to putin gain nothing cost ukrainian dearly jonah goldberg trump overture to putin gain nothing cost ukrainian dearly jonah goldberg trump overture to putin gain nothing cost ukrainian dearly jonah goldberg new electric vehicle groundbreaking automated sorter transform la crosse post office trump overture to putin gain nothing cost ukrainian dearly jonah goldberg trump overture to putin gain nothing cost ukrainian dearly jonah goldberg trump overture to putin gain nothing cost ukrainian dearly jonah goldberg trump overture to putin gain nothing cost ukrainian dearly jonah goldberg two reputed palestinian terrorist to speak at detroit conference obama applauds newsoms california redistricting plan a responsible a texas gop push new map obama applauds newsoms california redistricting plan a responsible a texas gop push new map obama applauds newsoms california redistricting plan a responsible a texas gop push new map obama applauds newsoms california redistricting plan a responsible a texas gop push new map texas house approves new map for congress requested by president trump military officer shifted to prosecute local dc crime amid trump takeover spectrum pharmaceutical sppi security lawsuit new lead plaintiff deadline announced parent of girl killed in texas camp flooding say common sense safety measure were absent a lawmaker consider bill newark public school launch newark save college saving program colorado secretary of state stomp on a projection rake while trying to pin a dictator label on trump ana navarro call out melania trump performative letter to putin ana navarro call out melania trump performative letter to putin trump team make mortgage fraud it newest revenge tool another judge denies doj request to release epstein grand jury document another judge denies doj request to release epstein grand jury document minnesota air national guard set to get super hercules fleet requiring new hangar socure win datos insight impact award in two fraud solution category best digital identity verification solution and best firstparty fraud solution route in newburgh closed because of tractortrailer collision judge block doj request to unseal grand jury transcript national advertising division find certain claim for kendamil infant formula supported the best toaster oven of tried and tested denver layoff deliver big hit to small office like childrens affair climate action county hoping to fill two committee our mental health in the hand of ai china is testing it first tiltrotor aircraft a man died gavin newscum and his ghoulish supporter share a laugh over trump assassination attempt a u navy transport ship just caught fire near okinawa wyoming lawmaker revisit ban on ballot drop box warren county to leave tourist rental rule in place bouloss family tie could help advance u national security interest in libya israel say it ha taken first step of military operation in gaza city vance derides protester while visiting national guard at union station vance derides protester while visiting national guard at union station cnn politics fed move to restrict funding for virginia school punishing support for transgender student court rule texas cannot mandate religious display in classroom court rule texas cannot mandate religious display in classroom new survey college student feel coerced into participating in animal experiment despite objection oppose such exercise epstein grand jury transcript wont be unsealed federal court rule newsomstyle redistricting effort critiqued by california democrat a recently a july statement show bill rhyne claim republican primary victory for spartanburg co sheriff sen mccormick of pennsylvania call to rebuild u military no act of lawfare is too petty for this administration voting poll for august court reject doj request to unseal jeffrey epstein grand jury file jd vance and pete hegseth visit national guard troop amid dc protest over trump crackdown u ramp up attack on international court over israel a look at those trump ha targeted in tactic of revoking security clearance'
$ws.Range("I42").Value = $headline
$ws.Range("J42").Value = $headline

# Sentiment + emotion numeric columns
$ws.Range("K42").Value = 0.179
$ws.Range("L42").Value = 0.125
$ws.Range("M42").Value = 0.695
$ws.Range("N42").Value = 0.9995000000000001
$ws.Range("O42").Value = 82
$ws.Range("P42").Value = 59
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = 161
$ws.Range("S42").Value = 76
$ws.Range("T42").Value = 198
$ws.Range("U42").Value = 118
$ws.Range("V42").Value = 58
$ws.Range("W42").Value = 19
$ws.Range("X42").Value = 67
